# Reorders the D/E/F/G columns of the SectorGroup codelist sheet.
#
# Old layout: D=group-code, E=category-name, F=group-name, G=category-code
# New layout: D=category-name, E=category-code, F=group-name, G=group-code
#
# i.e. new D = old E, new E = old G, new F = old F (unchanged), new G = old D.
# This is applied to every row (including the header row, whose labels follow
# the same rename/reorder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

$range = $ws.Range("A1").Resize($rowCount, $colCount)
$vals = $range.Value()

for ($r = 1; $r -le $rowCount; $r++) {
    $oldD = $vals[$r, 4]
    $oldE = $vals[$r, 5]
    $oldG = $vals[$r, 7]

    $vals[$r, 4] = $oldE
    $vals[$r, 5] = $oldG
    $vals[$r, 7] = $oldD
}

$range.Value = $vals
